$wb = $excel.ActiveWorkbook

# --- RegressionTests: fix B81 text and selection (whole row 1), before switching tabs ---
$reg = $wb.Worksheets.Item("RegressionTests")
$reg.Range("B81").Value = "Feature Integration: 1.7.2.3. Payload with push transition"
$reg.Range("A1:XFD1").Select() | Out-Null

# --- Insert new "Zinger" worksheet right after "RegressionTests" ---
$zinger = $wb.Worksheets.Add($null, $reg)
$zinger.Name = "Zinger"

# Column widths
$zinger.Columns.Item(1).ColumnWidth = 32
$zinger.Columns.Item(2).ColumnWidth = 45
$zinger.Columns.Item(6).ColumnWidth = 24
$zinger.Columns.Item(7).ColumnWidth = 43.140625
$zinger.Columns.Item(8).ColumnWidth = 33

# Header row formatting (bold, size 14, taller row)
$zinger.Range("A1:I1").Font.Bold = $true
$zinger.Range("A1:I1").Font.Size = 14
$zinger.Rows.Item(1).RowHeight = 18.75

# Header row values
$zinger.Range("A1").Value = "test"
$zinger.Range("B1").Value = "P feature tested "
$zinger.Range("C1").Value = "Static Error?"
$zinger.Range("D1").Value = "Dynamic Error?"
$zinger.Range("E1").Value = "Correct?"
$zinger.Range("F1").Value = "Runtime test?"
$zinger.Range("G1").Value = "Other features tested"
$zinger.Range("H1").Value = "Enabled error messages from P.4ml"
$zinger.Range("I1").Value = "Notes and TODOs"

# Row 3: BoundedChoice
$zinger.Range("A3").Value = "BoundedChoice"
$zinger.Range("B3").Value = "Zinger: testing bounded choice operator in Zing"
$zinger.Range("C3").Value = "No"
$zinger.Range("D3").Value = "No"
$zinger.Range("E3").Value = "Yes"
$zinger.Range("I3").Value = "Zinger arg: -bc:2"

# Row 4: BoundedChoice_1
$zinger.Range("A4").Value = "BoundedChoice_1"
$zinger.Range("B4").Value = "Zinger: testing bounded choice operator in Zing"
$zinger.Range("C4").Value = "No"
$zinger.Range("D4").Value = "Yes"
$zinger.Range("F4").Value = "Yes"
$zinger.Range("I4").Value = 'Same as BoundedChoice, but default "-bc" argument for zinger.exe'

# Row 5: DFSStackBound
$zinger.Range("A5").Value = "DFSStackBound"
$zinger.Range("B5").Value = "Zinger: testing DFS stack bound option in Zing"
$zinger.Range("C5").Value = "No"
$zinger.Range("D5").Value = "Yes"
$zinger.Range("F5").Value = "Yes"
$zinger.Range("H5").Value = '"DFS Stack Size Exceeded 10"'
$zinger.Range("I5").Value = 'Zinger args: "-maxdfsstack:10", "-ibound:100"'

# Select I5 as the active cell on the Zinger sheet (also makes Zinger the active/selected tab)
$zinger.Range("I5").Select() | Out-Null

Write-Output "done"
